# Added sheet number argument to processFile
# Replicates: new "Sheet2" (single new record) and new "Another Sheet"
# (all records, old + new) appended to the workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$headers = @("firstName", "lastName", "address.street", "address.city", "address.state", "address.zip")

# Re-select sheet1's original data range (matches the diff's updated selection)
$ws1.Range("A1:F3").Select()
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Sheet2: header row + the new "Max Irwin" record ---
$ws2 = $wb.Worksheets.Add($null, $ws1)

for ($c = 1; $c -le 6; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$ws2.Cells.Item(2, 1).Value = "Max"
$ws2.Cells.Item(2, 2).Value = "Irwin"
$ws2.Cells.Item(2, 3).Value = "123 Fake Street"
$ws2.Cells.Item(2, 4).Value = "Rochester"
$ws2.Cells.Item(2, 5).Value = "NY"
$ws2.Cells.Item(2, 6).Value = 99999

$ws2.Range("A2:F2").Select()
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- "Another Sheet": header row + all three records (existing two + new) ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Another Sheet"

for ($c = 1; $c -le 6; $c++) {
    $ws3.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$ws3.Cells.Item(2, 1).Value = "Jihad"
$ws3.Cells.Item(2, 2).Value = "Saladin"
$ws3.Cells.Item(2, 3).Value = "12 Beaver Court"
$ws3.Cells.Item(2, 4).Value = "Snowmass"
$ws3.Cells.Item(2, 5).Value = "CO"
$ws3.Cells.Item(2, 6).Value = 81615

$ws3.Cells.Item(3, 1).Value = "Marcus"
$ws3.Cells.Item(3, 2).Value = "Rivapoli"
$ws3.Cells.Item(3, 3).Value = "16 Vail Rd"
$ws3.Cells.Item(3, 4).Value = "Vail"
$ws3.Cells.Item(3, 5).Value = "CO"
$ws3.Cells.Item(3, 6).Value = 81657

$ws3.Cells.Item(4, 1).Value = "Max"
$ws3.Cells.Item(4, 2).Value = "Irwin"
$ws3.Cells.Item(4, 3).Value = "123 Fake Street"
$ws3.Cells.Item(4, 4).Value = "Rochester"
$ws3.Cells.Item(4, 5).Value = "NY"
$ws3.Cells.Item(4, 6).Value = 99999

# The first three rows (header + two pre-existing records) carry an explicit
# black font color (as if re-entered/pasted), the newly appended row does not.
$ws3.Range("A1:F3").Font.Color = 0

$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

$ws3.Range("F4").Select()

Write-Output "done"
